# Fill in Franziska Corradi's (row 6) availability for the week of
# columns L:Q ("A" = Abwesend / absent, "?" = Nicht sicher / unsure).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L6").Value = "A"
$ws.Range("M6").Value = "A"
$ws.Range("N6").Value = "?"
$ws.Range("O6").Value = "A"
$ws.Range("P6").Value = "?"
$ws.Range("Q6").Value = "?"

# Update the row's absence/unsure tally.
$ws.Range("AI6").Value = 6

# Move the active selection to reflect where the user last worked.
$ws.Range("Q6").Select()
